# Auto-generated: update FFXIV market price/profit data cells per scheduled runner refresh.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 591.25
$ws.Range("I6").Value = 165
$ws.Range("K6").Value = 495
$ws.Range("M6").Value = -383
$ws.Range("H9").Value = 63.153847
$ws.Range("I9").Value = 41.57143
$ws.Range("J9").Value = 88.333336
$ws.Range("K9").Value = 41.57143
$ws.Range("L9").Value = 88.333336
$ws.Range("M9").Value = 127.42857
$ws.Range("N9").Value = -426.333336
$ws.Range("H12").Value = 198.07692
$ws.Range("I12").Value = 206.25
$ws.Range("J12").Value = 100
$ws.Range("K12").Value = 206.25
$ws.Range("L12").Value = 100
$ws.Range("M12").Value = -36.25
$ws.Range("N12").Value = -440
$ws.Range("H21").Value = 18017
$ws.Range("J21").Value = 0
$ws.Range("L21").Value = 0
$ws.Range("N21").ClearContents()
$ws.Range("H23").Value = 18017
$ws.Range("J23").Value = 0
$ws.Range("L23").Value = 0
$ws.Range("N23").ClearContents()
$ws.Range("H28").Value = 585.1429000000001
$ws.Range("I28").Value = 603.5
$ws.Range("J28").Value = 475
$ws.Range("K28").Value = 603.5
$ws.Range("L28").Value = 475
$ws.Range("M28").Value = -118.5
$ws.Range("N28").Value = -1445
$ws.Range("H29").Value = 463.25
$ws.Range("I29").Value = 80.40000000000001
$ws.Range("J29").Value = 1101.3334
$ws.Range("K29").Value = 241.2
$ws.Range("L29").Value = 3304.0002
$ws.Range("M29").Value = 39.79999999999998
$ws.Range("N29").Value = -3866.0002
$ws.Range("H38").Value = 340.08334
$ws.Range("I38").Value = 198.27272
$ws.Range("J38").Value = 1900
$ws.Range("K38").Value = 594.81816
$ws.Range("L38").Value = 5700
$ws.Range("M38").Value = -222.81816
$ws.Range("N38").Value = -6444
$ws.Range("H58").Value = 1039.25
$ws.Range("I58").Value = 433.875
$ws.Range("J58").Value = 2250
$ws.Range("K58").Value = 1301.625
$ws.Range("L58").Value = 6750
$ws.Range("M58").Value = -1151.625
$ws.Range("N58").Value = -7050
$ws.Range("H74").Value = 3868.238
$ws.Range("I74").Value = 3842.7856
$ws.Range("J74").Value = 3919.1428
$ws.Range("K74").Value = 3842.7856
$ws.Range("L74").Value = 3919.1428
$ws.Range("M74").Value = -2906.7856
$ws.Range("N74").Value = -5791.1428
$ws.Range("H77").Value = 3868.238
$ws.Range("I77").Value = 3842.7856
$ws.Range("J77").Value = 3919.1428
$ws.Range("K77").Value = 19213.928
$ws.Range("L77").Value = 19595.714
$ws.Range("M77").Value = -14533.928
$ws.Range("N77").Value = -28955.714
$ws.Range("H111").Value = 1000
$ws.Range("I111").Value = 1000
$ws.Range("J111").Value = 0
$ws.Range("K111").Value = 3000
$ws.Range("L111").Value = 0
$ws.Range("M111").Value = 67
$ws.Range("N111").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 0
$ws.Range("I63").Value = 0
$ws.Range("K63").Value = 0
$ws.Range("M63").ClearContents()
$ws.Range("H66").Value = 0
$ws.Range("I66").Value = 0
$ws.Range("K66").Value = 0
$ws.Range("M66").ClearContents()
$ws.Range("H110").Value = 633.04
$ws.Range("I110").Value = 514.1739
$ws.Range("K110").Value = 514.1739
$ws.Range("M110").Value = 1530.8261
$ws.Range("H133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("L133").Value = 0

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H82").Value = 17410.75
$ws.Range("I82").Value = 8120
$ws.Range("J82").Value = 45283
$ws.Range("K82").Value = 8120
$ws.Range("L82").Value = 45283
$ws.Range("M82").Value = -7737
$ws.Range("N82").Value = -46049
$ws.Range("H85").Value = 17410.75
$ws.Range("I85").Value = 8120
$ws.Range("J85").Value = 45283
$ws.Range("K85").Value = 8120
$ws.Range("L85").Value = 45283
$ws.Range("M85").Value = -6794
$ws.Range("N85").Value = -47935
$ws.Range("H99").Value = 1044.36
$ws.Range("I99").Value = 1100.5
$ws.Range("J99").Value = 900
$ws.Range("K99").Value = 1100.5
$ws.Range("L99").Value = 900
$ws.Range("M99").Value = 397.5
$ws.Range("N99").Value = -3896
$ws.Range("H105").Value = 5122
$ws.Range("I105").Value = 6027.5
$ws.Range("J105").Value = 1500
$ws.Range("K105").Value = 6027.5
$ws.Range("L105").Value = 1500
$ws.Range("M105").Value = -4280.5
$ws.Range("N105").Value = -4994

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 3750
$ws.Range("I62").Value = 5000
$ws.Range("K62").Value = 5000
$ws.Range("M62").Value = -4376
$ws.Range("H65").Value = 3750
$ws.Range("I65").Value = 5000
$ws.Range("K65").Value = 25000
$ws.Range("M65").Value = -21880

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 25850.75
$ws.Range("J68").Value = 34334.332
$ws.Range("L68").Value = 103002.996
$ws.Range("N68").Value = -104624.996
$ws.Range("H71").Value = 25850.75
$ws.Range("J71").Value = 34334.332
$ws.Range("L71").Value = 309008.988
$ws.Range("N71").Value = -317120.988
$ws.Range("H131").Value = 847.7368
$ws.Range("J131").Value = 931.875
$ws.Range("L131").Value = 2795.625
$ws.Range("N131").Value = -12875.625

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 1642.6
$ws.Range("I113").Value = 0
$ws.Range("K113").Value = 0
$ws.Range("M113").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 378.05884
$ws.Range("I22").Value = 345.15384
$ws.Range("J22").Value = 485
$ws.Range("K22").Value = 345.15384
$ws.Range("L22").Value = 485
$ws.Range("M22").Value = -50.15384
$ws.Range("N22").Value = -1075
$ws.Range("H27").Value = 378.05884
$ws.Range("I27").Value = 345.15384
$ws.Range("J27").Value = 485
$ws.Range("K27").Value = 345.15384
$ws.Range("L27").Value = 485
$ws.Range("M27").Value = -238.15384
$ws.Range("N27").Value = -699

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 281.66666
$ws.Range("I107").Value = 279.53845
$ws.Range("J107").Value = 295.5
$ws.Range("K107").Value = 838.61535
$ws.Range("L107").Value = 886.5
$ws.Range("M107").Value = 1081.38465
$ws.Range("N107").Value = -4726.5
